# "Another Change has been made" -- add two more rows of data below the
# existing "1st change" cell, and touch up the sheet view the same way
# Excel would after typing into A2/A3 and reselecting/auto-sizing things.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data -- these land as new shared-string entries (index 1 and 2).
$ws.Range("A2").Value = "2nd change"
$ws.Range("A3").Value = "3rd change"

# Column A was best-fit to the new (wider) content.
$ws.Columns("A:A").ColumnWidth = 61/6

# Final cursor position left on D8 after the edits.
[void]$ws.Range("D8").Select()
